# Update the "Приступил" time entry in row 33 of sheet "Лист1":
#  - C33: "Приступил  8:30/11:00" -> "Приступил  12:00/13:00"
#  - E33: 0 -> 1 (hours worked)
# E35 is a SUM formula over E18:E34 and will recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("C33").Value = "Приступил  12:00/13:00"
$ws.Range("E33").Value = 1

# Move the saved selection/scroll position to F33 (matches the final
# on-screen state after making the edit), scrolled back to the top.
$ws.Activate()
$ws.Range("F33").Select()
$excel.ActiveWindow.ScrollRow = 1
